# Swap the values of columns C and D (codeforiati:group-name and
# codeforiati:group-code, including the header row) for every used row
# on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
